$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 28
$ws.Range("I11").Value = 28
$ws.Range("K11").Value = 28
$ws.Range("M11").Value = 112

$ws.Range("H17").Value = 1759.8
$ws.Range("J17").Value = 1712.25
$ws.Range("L17").Value = 5136.75
$ws.Range("N17").Value = -5472.75

$ws.Range("H32").Value = 8166.6665
$ws.Range("J32").Value = 8166.6665
$ws.Range("L32").Value = 8166.6665
$ws.Range("N32").Value = -8818.666499999999

$ws.Range("H92").Value = 2854.2856
$ws.Range("I92").Value = 2747
$ws.Range("J92").Value = 2997.3333
$ws.Range("K92").Value = 2747
$ws.Range("L92").Value = 2997.3333
$ws.Range("M92").Value = -1499
$ws.Range("N92").Value = -5493.3333

$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()

$ws.Range("H135").Value = 1962.6364
$ws.Range("I135").Value = 1911.2222
$ws.Range("K135").Value = 17200.9998
$ws.Range("M135").Value = -14665.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2361.7273
$ws.Range("I2").Value = 2613.7144
$ws.Range("J2").Value = 1920.75
$ws.Range("K2").Value = 2613.7144
$ws.Range("L2").Value = 1920.75
$ws.Range("M2").Value = -2500.7144
$ws.Range("N2").Value = -2146.75

$ws.Range("H45").Value = 2764.2856
$ws.Range("I45").Value = 2270
$ws.Range("J45").Value = 4000
$ws.Range("K45").Value = 2270
$ws.Range("L45").Value = 4000
$ws.Range("M45").Value = -1893
$ws.Range("N45").Value = -4754

$ws.Range("H116").Value = 2361.7273
$ws.Range("I116").Value = 2613.7144
$ws.Range("J116").Value = 1920.75
$ws.Range("K116").Value = 2613.7144
$ws.Range("L116").Value = 1920.75
$ws.Range("M116").Value = -319.7143999999998
$ws.Range("N116").Value = -6508.75

$ws.Range("H122").Value = 7113.5557
$ws.Range("I122").Value = 8080.846
$ws.Range("K122").Value = 24242.538
$ws.Range("M122").Value = -21792.538

$ws.Range("H132").Value = 3282
$ws.Range("I132").Value = 3329
$ws.Range("K132").Value = 9987
$ws.Range("M132").Value = -7457

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2361.7273
$ws.Range("I3").Value = 2613.7144
$ws.Range("J3").Value = 1920.75
$ws.Range("K3").Value = 2613.7144
$ws.Range("L3").Value = 1920.75
$ws.Range("M3").Value = -2499.7144
$ws.Range("N3").Value = -2148.75

$ws.Range("H20").Value = 5000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 5000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 5000
$ws.Range("N20").Value = -5494
$ws.Range("M20").ClearContents()

$ws.Range("H80").Value = 751.375
$ws.Range("I80").Value = 584.6667
$ws.Range("K80").Value = 584.6667
$ws.Range("M80").Value = 413.3333

$ws.Range("H82").Value = 16013.8
$ws.Range("I82").Value = 16013.8
$ws.Range("K82").Value = 16013.8
$ws.Range("M82").Value = -15630.8

$ws.Range("H83").Value = 751.375
$ws.Range("I83").Value = 584.6667
$ws.Range("K83").Value = 2923.3335
$ws.Range("M83").Value = 2068.6665

$ws.Range("H85").Value = 16013.8
$ws.Range("I85").Value = 16013.8
$ws.Range("K85").Value = 16013.8
$ws.Range("M85").Value = -14687.8

$ws.Range("H97").Value = 12432.667
$ws.Range("I97").Value = 12432.667
$ws.Range("K97").Value = 12432.667
$ws.Range("M97").Value = -11441.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 126.4
$ws.Range("I7").Value = 126.4
$ws.Range("K7").Value = 126.4
$ws.Range("M7").Value = -13.40000000000001

$ws.Range("H22").Value = 1250
$ws.Range("I22").Value = 1250
$ws.Range("K22").Value = 1250
$ws.Range("M22").Value = -900

$ws.Range("H35").Value = 913
$ws.Range("I35").Value = 913
$ws.Range("K35").Value = 913
$ws.Range("M35").Value = -619

$ws.Range("H36").Value = 548
$ws.Range("I36").Value = 548
$ws.Range("K36").Value = 548
$ws.Range("M36").Value = -160

$ws.Range("H40").Value = 548
$ws.Range("I40").Value = 548
$ws.Range("K40").Value = 548
$ws.Range("M40").Value = -388

$ws.Range("H41").Value = 12439.8
$ws.Range("J41").Value = 20000
$ws.Range("L41").Value = 20000
$ws.Range("N41").Value = -20856

$ws.Range("H47").Value = 2000
$ws.Range("I47").Value = 2000
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 2000
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -1434
$ws.Range("N47").ClearContents()

$ws.Range("H48").Value = 1200
$ws.Range("I48").Value = 1200
$ws.Range("K48").Value = 1200
$ws.Range("M48").Value = -724

$ws.Range("H99").Value = 1200
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws.Range("H126").Value = 1200
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws.Range("H132").Value = 2222.2
$ws.Range("I132").Value = 2222.2
$ws.Range("K132").Value = 6666.599999999999
$ws.Range("M132").Value = -4136.599999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 874.625
$ws.Range("I23").Value = 819
$ws.Range("J23").Value = 967.3333
$ws.Range("K23").Value = 2457
$ws.Range("L23").Value = 2901.9999
$ws.Range("M23").Value = -2222
$ws.Range("N23").Value = -3371.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 17002666
$ws.Range("I7").Value = 25500000
$ws.Range("J7").Value = 8000
$ws.Range("K7").Value = 25500000
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = -25499888
$ws.Range("N7").Value = -8224

$ws.Range("H8").Value = 17002666
$ws.Range("I8").Value = 25500000
$ws.Range("J8").Value = 8000
$ws.Range("K8").Value = 25500000
$ws.Range("L8").Value = 8000
$ws.Range("M8").Value = -25499861
$ws.Range("N8").Value = -8278

$ws.Range("H122").Value = 6785.2856
$ws.Range("I122").Value = 5499.6665
$ws.Range("K122").Value = 16498.9995
$ws.Range("M122").Value = -14048.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 6604.5
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H16").Value = 5000
$ws.Range("I16").Value = 5000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 5000
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -4830
$ws.Range("N16").ClearContents()

$ws.Range("H20").Value = 9002.5
$ws.Range("I20").Value = 8005
$ws.Range("K20").Value = 8005
$ws.Range("M20").Value = -7779

$ws.Range("H22").Value = 733.1111
$ws.Range("I22").Value = 762.25
$ws.Range("K22").Value = 762.25
$ws.Range("M22").Value = -467.25

$ws.Range("H27").Value = 733.1111
$ws.Range("I27").Value = 762.25
$ws.Range("K27").Value = 762.25
$ws.Range("M27").Value = -655.25

$ws.Range("H40").Value = 4276.6
$ws.Range("J40").Value = 3944
$ws.Range("L40").Value = 3944
$ws.Range("N40").Value = -4216

$ws.Range("H55").Value = 3400
$ws.Range("I55").Value = 800
$ws.Range("K55").Value = 800
$ws.Range("M55").Value = -627

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 7500000
$ws.Range("I21").Value = 7500000
$ws.Range("K21").Value = 7500000
$ws.Range("M21").Value = -7499765

$ws.Range("H24").Value = 5000000
$ws.Range("I24").Value = 5000000
$ws.Range("K24").Value = 5000000
$ws.Range("M24").Value = -4999770

$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()

$ws.Range("H30").Value = 40000
$ws.Range("I30").Value = 40000
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 40000
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -39893
$ws.Range("N30").ClearContents()

$ws.Range("H35").Value = 7500000
$ws.Range("I35").Value = 7500000
$ws.Range("K35").Value = 7500000
$ws.Range("M35").Value = -7499710

$ws.Range("H122").Value = 1728.5
$ws.Range("I122").Value = 683.6
$ws.Range("K122").Value = 2050.8
$ws.Range("M122").Value = 399.1999999999998

Write-Output "applied all changes"
